$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1, matching the formatting of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" column values for rows 2..8 (plain numbers, no special style)
$saveValues = @(0, 0, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
